# Handles float input without breaking stuff
# Updates the marksheet summary (rows 10-12) and the per-question
# "Student Ans" grid so it reflects a 28-question quiz (instead of the
# original 56-question / all-blank layout), and drops the unused third
# "Student Ans / Correct Ans" block (columns G:H) plus most of the
# second block (columns D:E), which are no longer part of the quiz.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Summary block (rows 10-12)
# ---------------------------------------------------------------

# Row 10 ("No." row: Right / Wrong / Not Attempt / Max)
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 28

# Row 11 ("Marking" row: marks per right / wrong answer)
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 ("Total" row: total marks + final score)
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -5
$ws.Range("E12").Value = "63/112"

# ---------------------------------------------------------------
# Drop the third "Student Ans / Correct Ans" block (columns G:H) --
# it no longer corresponds to any question.
# ---------------------------------------------------------------
$ws.Range("G15:H40").Clear()

# ---------------------------------------------------------------
# Second "Student Ans / Correct Ans" block (columns D:E) now only
# covers 3 questions (rows 16-18); clear the rest (rows 19-40).
# ---------------------------------------------------------------
$ws.Range("D19:E40").Clear()

$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D17").Value = "Option B"
$ws.Range("D17").Style = "incorrectStyle"
$ws.Range("D18").Style = "normalStyle"

# ---------------------------------------------------------------
# First "Student Ans" block (column A) now holds the student's
# actual answers (correct/incorrect/blank), instead of being left
# empty on every row.
# ---------------------------------------------------------------
$ws.Range("A16").Value = "Option A"
$ws.Range("A16").Style = "correctStyle"
$ws.Range("A17").Value = "Option A"
$ws.Range("A17").Style = "incorrectStyle"
$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"
$ws.Range("A20").Value = "Option B"
$ws.Range("A20").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"
$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"
$ws.Range("A23").Value = "Option D"
$ws.Range("A23").Style = "correctStyle"
$ws.Range("A25").Value = "Option A"
$ws.Range("A25").Style = "correctStyle"
$ws.Range("A26").Value = "Option C"
$ws.Range("A26").Style = "correctStyle"
$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"
$ws.Range("A28").Value = "Option B"
$ws.Range("A28").Style = "incorrectStyle"
$ws.Range("A29").Value = "Option D"
$ws.Range("A29").Style = "correctStyle"
$ws.Range("A30").Value = "Option B"
$ws.Range("A30").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"
$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"
$ws.Range("A34").Value = "Option B"
$ws.Range("A34").Style = "correctStyle"
$ws.Range("A36").Value = "Option D"
$ws.Range("A36").Style = "incorrectStyle"
$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"
$ws.Range("A40").Value = "Option B"
$ws.Range("A40").Style = "incorrectStyle"
# A19, A24, A31, A35, A37 stay blank ("normalStyle") -- not attempted.
